$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.6
$ws.Range("K2").Value = 1.77
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 2.45
$ws.Range("K3").Value = 1.73
$ws.Range("N3").Value = 4.75
$ws.Range("S3").Value = 1.77
$ws.Range("T3").Value = 2
$ws.Range("AC3").Value = 4.75
$ws.Range("AQ3").Value = 101
$ws.Range("AT3").Value = 1.91
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 1.83
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 1.62
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("Y4").Value = 11
$ws.Range("AA4").Value = 26
$ws.Range("AC4").Value = 6
$ws.Range("AE4").Value = 21
$ws.Range("AK4").Value = 34
$ws.Range("AP4").Value = 34
$ws.Range("AR4").Value = 101
$ws.Range("AS4").Value = 351
$ws.Range("AT4").Value = 2.2
$ws.Range("AU4").Value = 9.5
$ws.Range("AX4").Value = 21
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 401
$ws.Range("G6").Value = 2.3
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.85
$ws.Range("G7").Value = 1.68
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 4.5
$ws.Range("J7").Value = 2.3
$ws.Range("K7").Value = 2.3
$ws.Range("L7").Value = 4.75
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.08
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25
$ws.Range("W7").Value = 8.5
$ws.Range("X7").Value = 9
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 13
$ws.Range("AD7").Value = 7.5
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 15
$ws.Range("AJ7").Value = 51
$ws.Range("AK7").Value = 34
$ws.Range("AM7").Value = 151
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 9
$ws.Range("AR7").Value = 41
$ws.Range("AT7").Value = 3.25
$ws.Range("AU7").Value = 8
$ws.Range("AW7").Value = 6.5
$ws.Range("AX7").Value = 23
$ws.Range("AY7").Value = 29
$ws.Range("AZ7").Value = 81
$ws.Range("G8").Value = 1.71
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("G9").Value = 1.62
$ws.Range("G10").Value = 2.82
$ws.Range("I10").Value = 2.55
$ws.Range("G11").Value = 2.45
$ws.Range("I11").Value = 2.75
$ws.Range("G12").Value = 1.81
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.4
$ws.Range("Q12").Value = 2.07
$ws.Range("R12").Value = 1.83
$ws.Range("I14").Value = 3.5
$ws.Range("U14").Value = 1.83
$ws.Range("V14").Value = 1.83
$ws.Range("W14").Value = 7
$ws.Range("X14").Value = 9.5
$ws.Range("Z14").Value = 19
$ws.Range("AA14").Value = 19
$ws.Range("AG14").Value = 9
$ws.Range("AI14").Value = 12
$ws.Range("AO14").Value = 12
$ws.Range("N15").Value = 10
$ws.Range("G17").Value = 1.42
$ws.Range("Y17").Value = 9
$ws.Range("AE17").Value = 21
$ws.Range("AM17").Value = 401
$ws.Range("M18").Value = 1.07
$ws.Range("O18").Value = 1.36
$ws.Range("M19").Value = 1.06
$ws.Range("O19").Value = 1.3
$ws.Range("Q19").Value = 1.92
$ws.Range("R19").Value = 1.82
$ws.Range("M20").Value = 1.08
$ws.Range("O20").Value = 1.4
$ws.Range("M21").Value = 1.03
$ws.Range("O21").Value = 1.18
$ws.Range("M22").Value = 1.04
$ws.Range("O22").Value = 1.22
$ws.Range("M23").Value = 1.06
$ws.Range("O23").Value = 1.29
$ws.Range("Q23").Value = 1.94
$ws.Range("R23").Value = 1.79
$ws.Range("O24").Value = 1.29
$ws.Range("P24").Value = 3.5
$ws.Range("Q24").Value = 1.87
$ws.Range("R24").Value = 1.87
$ws.Range("Q25").Value = 1.82
$ws.Range("R25").Value = 1.92
$ws.Range("Q26").Value = 2.1
$ws.Range("R26").Value = 1.67
$ws.Range("G28").Value = 1.95
$ws.Range("H28").Value = 3.3
$ws.Range("I28").Value = 3.9
$ws.Range("J28").Value = 2.63
$ws.Range("K28").Value = 2.1
$ws.Range("M28").Value = 1.07
$ws.Range("N28").Value = 9
$ws.Range("Q28").Value = 2.1
$ws.Range("R28").Value = 1.67
$ws.Range("U28").Value = 1.91
$ws.Range("V28").Value = 1.91
$ws.Range("AC28").Value = 9
$ws.Range("AD28").Value = 6.5
$ws.Range("AR28").Value = 51
$ws.Range("AZ28").Value = 81
$ws.Range("U29").Value = 1.75
$ws.Range("U31").Value = 1.57
$ws.Range("V31").Value = 2.25
$ws.Range("Z31").Value = 13
$ws.Range("AG31").Value = 17
$ws.Range("AH31").Value = 26
$ws.Range("AX31").Value = 23
$ws.Range("AY31").Value = 26
$ws.Range("BB31").Value = 151
$ws.Range("AE32").Value = 15
$ws.Range("AO32").Value = 26
$ws.Range("AW32").Value = 3.75
$ws.Range("G33").Value = 1.53
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = 3.25
$ws.Range("I34").Value = 2.3
$ws.Range("J34").Value = 3.5
$ws.Range("L34").Value = 3
$ws.Range("M34").Value = 1.05
$ws.Range("N34").Value = 11
$ws.Range("AA34").Value = 23
$ws.Range("AH34").Value = 12
$ws.Range("AI34").Value = 9.5
$ws.Range("AM34").Value = 201
$ws.Range("AN34").Value = 5
$ws.Range("AO34").Value = 17
$ws.Range("AY34").Value = 21
$ws.Range("BA34").Value = 51
$ws.Range("M35").Value = 1.04
$ws.Range("O35").Value = 1.22
$ws.Range("S35").Value = 1.3
$ws.Range("M36").Value = 1.07
$ws.Range("O36").Value = 1.36
$ws.Range("S36").Value = 1.47
$ws.Range("M37").Value = 1.05
$ws.Range("O37").Value = 1.29
$ws.Range("M38").Value = 1.05
$ws.Range("O38").Value = 1.29
$ws.Range("V39").Value = 1.57
